$d = $word.ActiveDocument

# 1) Text change: the venue-name placeholder is replaced by the external
#    short-name placeholder, and the stray space after <<else>> is removed.
$d.Content.Find.Execute(
    ".venue_name>><<else>> Online Civil Claims<<es_>>",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ".external_short_name>><<else>>Online Civil Claims<<es_>>",
    2) | Out-Null

# 2) Formatting change: the heading paragraph ("In the County Court at...")
#    gains line spacing of 1.15 (line=276, auto rule).
$rng = $d.Content
$rng.Find.Execute("In the County Court at", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $rng.Paragraphs(1)
$para.LineSpacingRule = 5   # wdLineSpaceMultiple
$para.LineSpacing = 13.8    # 276 twentieths-of-a-point / 20 = 13.8 points
